$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 91: rolling 7-day count / rate revised upward by one case -----
$ws.Range("C91").Value = 6
$ws.Range("D91").Value = 97.0402717127608

# --- Insert a brand-new data row (date 44235) right after row 92 -------
# This shifts the former rows 93-113 down to 94-114.
$ws.Rows.Item(93).Insert()

# Copy formatting (date number format, borders, font, etc.) from the row
# above down onto the freshly inserted row so it matches its neighbours.
$ws.Range("A92:D92").Copy()
$ws.Range("A93:D93").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A93").Value = 44235
$ws.Range("B93").Value = 1
$ws.Range("C93").Value = 6
$ws.Range("D93").Value = 97.0402717127608

# --- A few rows after the insertion need their rolling-window values ---
# recomputed (dates/new-case counts already shifted correctly by the
# row insertion above; only C/D need new totals).
$ws.Range("C94").Value = 11
$ws.Range("D94").Value = 177.9071648067281

$ws.Range("C95").Value = 13
$ws.Range("D95").Value = 210.2539220443151

$ws.Range("C96").Value = 19
$ws.Range("D96").Value = 307.2941937570758

# --- Row 112 (shifted from the old row 111) now gets rolling-window ----
# totals it previously lacked.
$ws.Range("C112").Value = 22
$ws.Range("D112").Value = 355.8143296134562

# --- Append one brand-new row (date 44257) at the bottom of the table --
# Copy formatting from the previous last row down onto the new one, then
# set its date and new-case count (C/D rolling-window data is not yet
# available for this day, so those cells are left blank).
$ws.Range("A114:D114").Copy()
$ws.Range("A115:D115").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A115").Value = 44257
$ws.Range("B115").Value = 6
